# Enjoy-SQL.pptx edit:
#   Slide 4 contains a SmartArt (radial) diagram. The "Assumption" node's
#   text currently reads  "Assumption:Player can´t SQL"  and needs to
#   become "Assumption: Player doesn't know SQL".
#
# We reach the diagram through the normal PowerPoint SmartArt object
# model: Shape.SmartArt -> SmartArt.AllNodes(i) -> Node.TextFrame2.TextRange.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# "Inhaltsplatzhalter 8" is the content placeholder holding the SmartArt graphic.
$sh = $s.Shapes.Item(2)

$sa = $sh.SmartArt
$nodes = $sa.AllNodes

for ($i = 1; $i -le $nodes.Count; $i++) {
    $node = $nodes.Item($i)
    $txt = $node.TextFrame2.TextRange.Text
    if ($txt -like "Assumption*SQL") {
        $node.TextFrame2.TextRange.Text = "Assumption: Player doesn't know SQL"
    }
}

Write-Output "done"
